$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Through 2022-12-02" to "Through 2022-12-03"
$ws.Name = "Through 2022-12-03"

# Update header label in I1 from "2022 (through 12-02)" to "2022 (through 12-03)"
$ws.Range("I1").Value = "2022 (through 12-03)"

# Update the December 2022 data: I13 (December) 10 -> 13, I14 (Total) 1526 -> 1529
$ws.Range("I13").Value = 13
$ws.Range("I14").Value = 1529
